# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts.
#
# This updates the DAMSLTag (column I) and DialogAct (column J) values
# for a set of rows in Sheet1, reflecting the re-annotated dialog acts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> (DAMSLTag, DialogAct)
$updates = @{
    2   = @("sv", "Statement-opinion")
    5   = @("b", "Acknowledge (Backchannel)")
    7   = @("sv", "Statement-opinion")
    8   = @("sv", "Statement-opinion")
    13  = @("b", "Acknowledge (Backchannel)")
    19  = @("%", "Uninterpretable")
    30  = @("sd", "Statement-non-opinion")
    35  = @("aa", "Agree/Accept")
    48  = @("aa", "Agree/Accept")
    52  = @("aa", "Agree/Accept")
    57  = @("sv", "Statement-opinion")
    58  = @("%", "Uninterpretable")
    62  = @("sd", "Statement-non-opinion")
    66  = @("sv", "Statement-opinion")
    68  = @("sd", "Statement-non-opinion")
    79  = @("sd", "Statement-non-opinion")
    82  = @("sd", "Statement-non-opinion")
    83  = @("%", "Uninterpretable")
    96  = @("aa", "Agree/Accept")
    97  = @("sd", "Statement-non-opinion")
    100 = @("ba", "Appreciation")
    102 = @("aa", "Agree/Accept")
    114 = @("sd", "Statement-non-opinion")
    119 = @("sd", "Statement-non-opinion")
    129 = @("sd", "Statement-non-opinion")
    134 = @("aa", "Agree/Accept")
    137 = @("sv", "Statement-opinion")
    150 = @("sd", "Statement-non-opinion")
    170 = @("aa", "Agree/Accept")
    179 = @("sv", "Statement-opinion")
    183 = @("sd", "Statement-non-opinion")
    184 = @("sd", "Statement-non-opinion")
    185 = @("sd", "Statement-non-opinion")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
